$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the close value of the last existing row (row 235)
$ws.Range("F235").Value2 = 457.82

# 2. Append three new monthly rows (236-238) of FX data, matching the
#    existing layout: A=datetime (serial), B=symbol, C=open, D=high,
#    E=low, F=close, G=volume.

$newRows = @(
    @{ Row = 236; DateTime = 45170.33333333334; Open = 457.82; High = 485.93; Low = 454.52; Close = 477.37; Volume = 0 },
    @{ Row = 237; DateTime = 45201.375;          Open = 477.37; High = 480.77; Low = 467.91; Close = 468.27; Volume = 0 },
    @{ Row = 238; DateTime = 45231.375;          Open = 468.27; High = 470.91; Low = 461.68; Close = 464.87; Volume = 0 }
)

foreach ($r in $newRows) {
    $rowIndex = $r.Row

    # Copy the date cell format (numeric date style with border/alignment)
    # from the previous row's column A cell so the new cell matches the
    # existing style used throughout column A.
    $ws.Range("A" + ($rowIndex - 1)).Copy()
    $ws.Range("A" + $rowIndex).PasteSpecial(-4122)  # xlPasteFormats
    $excel.CutCopyMode = 0

    $ws.Cells.Item($rowIndex, 1).Value2 = $r.DateTime
    $ws.Cells.Item($rowIndex, 2).Value2 = "FX_IDC:USDKZT"
    $ws.Cells.Item($rowIndex, 3).Value2 = $r.Open
    $ws.Cells.Item($rowIndex, 4).Value2 = $r.High
    $ws.Cells.Item($rowIndex, 5).Value2 = $r.Low
    $ws.Cells.Item($rowIndex, 6).Value2 = $r.Close
    $ws.Cells.Item($rowIndex, 7).Value2 = $r.Volume
}
